$d = $word.ActiveDocument

# Locate the paragraph that ends the "uparei no GitHub." sentence - it is
# the last non-empty paragraph in the document (just before the trailing
# empty paragraph / sectPr).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*uparei no GitHub*") {
        $targetIndex = $i
    }
}

$anchor = $d.Paragraphs.Item($targetIndex).Range
$anchor.Collapse(0)

$rpr = "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>"
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 1: the "-------- atualização do relatório --------" divider
[void]$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Item($targetIndex + 1).Range
$xml1 = "<w:p $wns><w:pPr>$rpr</w:pPr>" + `
        "<w:r>$rpr<w:t>-------------------------------------- atualização do relatório</w:t></w:r>" + `
        "<w:r>$rpr<w:t>--------------------------------------</w:t></w:r>" + `
        "</w:p>"
[void]$p.InsertXML($xml1)

# --- Paragraph 2: empty paragraph
$p2anchor = $d.Paragraphs.Item($targetIndex + 1).Range
$p2anchor.Collapse(0)
[void]$p2anchor.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($targetIndex + 2).Range
$xml2 = "<w:p $wns><w:pPr>$rpr</w:pPr></w:p>"
[void]$p2.InsertXML($xml2)

# --- Paragraph 3: "Eu programei e organizei ..." with proofErr spans
$p3anchor = $d.Paragraphs.Item($targetIndex + 2).Range
$p3anchor.Collapse(0)
[void]$p3anchor.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($targetIndex + 3).Range

$xml3 = "<w:p $wns><w:pPr>$rpr</w:pPr>" + `
        "<w:r>$rpr<w:t xml:space='preserve'>Eu programei e organizei todos os </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r>$rpr<w:t>Assets</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r>$rpr<w:t xml:space='preserve'> na </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r>$rpr<w:t>Unity</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r>$rpr<w:t xml:space='preserve'>, ajudei na coordenação do projeto junto ao Gustavo </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r>$rpr<w:t>Hanada</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r>$rpr<w:t xml:space='preserve'>, ajudei na gravação dos áudios utilizados no projeto, auxiliei o grupo na utilização do </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r>$rpr<w:t>Github</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r>$rpr<w:t xml:space='preserve'>. </w:t></w:r>" + `
        "</w:p>"
[void]$p3.InsertXML($xml3)
